$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) from Excel auto-converting numeric-looking
# text (e.g. "221.21", "30.719.47") into actual numbers, by temporarily
# forcing a Text number format while assigning the values, then clearing
# the format again so the cell style matches the original (unstyled) cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.719.47"
$ws.Range("D3").Value = "1.689.61"
$ws.Range("D5").Value = "221.21"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "30.70"
$ws.Range("D9").Value = "0.265"
$ws.Range("D12").Value = "1.933.58"
$ws.Range("D13").Value = "10.65"
$ws.Range("D14").Value = "0.623"
$ws.Range("D15").Value = "1.694.93"
$ws.Range("D16").Value = "3.99"
$ws.Range("D17").Value = "30.736.08"
$ws.Range("D18").Value = "66.46"
$ws.Range("D19").Value = "247.04"
$ws.Range("D20").Value = "0.0₃0716"
$ws.Range("D22").Value = "10.28"
$ws.Range("D25").Value = "157.20"
$ws.Range("D26").Value = "15.90"
$ws.Range("D28").Value = "6.72"
$ws.Range("D33").Value = "1.515.48"
$ws.Range("D37").Value = "83.62"
$ws.Range("D40").Value = "2.71"
$ws.Range("D42").Value = "0.849"
$ws.Range("D43").Value = "0.0505"
$ws.Range("D44").Value = "2.01"
$ws.Range("D45").Value = "1.04"
$ws.Range("D46").Value = "1.00"
$ws.Range("D47").Value = "51.91"
$ws.Range("D48").Value = "1.825.98"
$ws.Range("D49").Value = "5.45"
$ws.Range("D50").Value = "95.09"
$ws.Range("D51").Value = "0.0₆0116"

$ws.Range("D2:D51").ClearFormats()

# Remaining columns (Coin name, Link, Volume%) are plain text already and
# do not get reinterpreted as numbers by Excel, so they can be set directly.
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E13").Value = "  +12.32%  "
$ws.Range("E14").Value = "  +8.23%  "
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +7.58%  "
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("E39").Value = "  +4.51%  "
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -6.77%  "
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  +5.03%  "
$ws.Range("E51").Value = "  +1.74%  "
